$p = $ppt.ActivePresentation

# Delete the "Class Activity" slide (position 18) and the short
# "Next Time on PSY 203..." slide (position 19). Delete from the
# higher index first so the lower index stays valid.
$p.Slides.Item(19).Delete()
$p.Slides.Item(18).Delete()

# Move "Measuring Correlation" (now at position 26, after the two
# deletions) up to position 17, right after "Correlation at a Glance"
# and before "Limitations of a Correlation Study".
$p.Slides.Item(26).MoveTo(17)
